$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text by pre-formatting as Text ("@")
# for the specific cells whose new value could otherwise be auto-converted to a number.
$textCells = @(
    'D4',
    'D5',
    'D6',
    'D8',
    'D9',
    'D10',
    'D14',
    'D18',
    'D19',
    'D20',
    'D21',
    'D22',
    'D23',
    'D25',
    'D27',
    'D28',
    'D31',
    'D32',
    'D33',
    'D34',
    'D35',
    'D37',
    'D39',
    'D40',
    'D41',
    'D43',
    'D44',
    'D45',
    'D46',
    'D47',
    'D48',
    'D49',
    'D51'
)
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated values (prices, volumes, and a few re-ranked coin name/link swaps)
$ws.Range('D2').Value = '70.137.24'
$ws.Range('E2').Value = '  +0.88%  '
$ws.Range('D3').Value = '3.708.11'
$ws.Range('E3').Value = '  +0.44%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '651.20'
$ws.Range('E5').Value = '  -3.92%  '
$ws.Range('D6').Value = '162.86'
$ws.Range('E6').Value = '  +0.75%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '0.505'
$ws.Range('E8').Value = '  +2.07%  '
$ws.Range('D9').Value = '0.148'
$ws.Range('E9').Value = '  +0.60%  '
$ws.Range('D10').Value = '7.20'
$ws.Range('E10').Value = '  +1.00%  '
$ws.Range('E11').Value = '  +1.89%  '
$ws.Range('E12').Value = '  +0.96%  '
$ws.Range('D13').Value = '4.325.04'
$ws.Range('E13').Value = '  +0.24%  '
$ws.Range('D14').Value = '33.07'
$ws.Range('E14').Value = '  +1.79%  '
$ws.Range('D15').Value = '3.700.82'
$ws.Range('E15').Value = '  -0.42%  '
$ws.Range('D16').Value = '70.080.04'
$ws.Range('E16').Value = '  +0.88%  '
$ws.Range('E17').Value = '  +0.54%  '
$ws.Range('D18').Value = '16.23'
$ws.Range('E18').Value = '  +1.39%  '
$ws.Range('D19').Value = '6.57'
$ws.Range('E19').Value = '  +1.58%  '
$ws.Range('D20').Value = '10.65'
$ws.Range('E20').Value = '  +8.86%  '
$ws.Range('D21').Value = '475.73'
$ws.Range('E21').Value = '  +1.10%  '
$ws.Range('D22').Value = '0.655'
$ws.Range('E22').Value = '  +0.76%  '
$ws.Range('D23').Value = '80.33'
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('D24').Value = '3.850.94'
$ws.Range('E24').Value = '  +0.34%  '
$ws.Range('D25').Value = '0.0000130'
$ws.Range('E25').Value = '  +3.37%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').Value = '11.11'
$ws.Range('E27').Value = '  +2.32%  '
$ws.Range('D28').Value = '9.30'
$ws.Range('E28').Value = '  +1.75%  '
$ws.Range('E29').Value = '  -1.02%  '
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('D31').Value = '2.03'
$ws.Range('E31').Value = '  +0.86%  '
$ws.Range('D32').Value = '6.63'
$ws.Range('E32').Value = '  +0.67%  '
$ws.Range('D33').Value = '27.09'
$ws.Range('E33').Value = '  +0.30%  '
$ws.Range('B34').Value = 'Kaspa'
$ws.Range('C34').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D34').Value = '0.167'
$ws.Range('E34').Value = '  +2.99%  '
$ws.Range('B35').Value = 'Binance-PegBSC-USD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.67%  '
$ws.Range('D36').Value = '3.700.79'
$ws.Range('E36').Value = '  +0.51%  '
$ws.Range('D37').Value = '8.52'
$ws.Range('E37').Value = '  +0.30%  '
$ws.Range('E38').Value = '  -0.07%  '
$ws.Range('D39').Value = '2.30'
$ws.Range('E39').Value = '  +1.29%  '
$ws.Range('D40').Value = '5.95'
$ws.Range('E40').Value = '  -4.22%  '
$ws.Range('D41').Value = '180.60'
$ws.Range('E41').Value = '  +7.14%  '
$ws.Range('E42').Value = '  -0.10%  '
$ws.Range('D43').Value = '0.0912'
$ws.Range('E43').Value = '  +1.29%  '
$ws.Range('D44').Value = '0.934'
$ws.Range('E44').Value = '  -0.91%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').Value = '2.89'
$ws.Range('E45').Value = '  +5.47%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').Value = '47.15'
$ws.Range('E46').Value = '  +1.00%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = '29.41'
$ws.Range('E47').Value = '  +6.00%  '
$ws.Range('D48').Value = '0.000277'
$ws.Range('E48').Value = '  -0.38%  '
$ws.Range('D49').Value = '1.28'
$ws.Range('E49').Value = '  -1.60%  '
$ws.Range('E50').Value = '  -1.54%  '
$ws.Range('D51').Value = '7.90'
$ws.Range('E51').Value = '  +0.14%  '
